# Update the "last saved" date field text shown in the slide master and
# every slide layout's date placeholder (2023/8/14 -> 2023/12/6), and
# update the title text on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "2023/8/14"
$newDate = "2023/12/6"

function Update-DateShapeText($shapes) {
    $updated = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
                $updated = $updated + 1
            }
        }
    }
    return $updated
}

# Slide master footer date placeholder.
$sm = $p.SlideMaster
Update-DateShapeText $sm.Shapes | Out-Null

# Every slide layout's footer date placeholder.
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    Update-DateShapeText $cl.Shapes | Out-Null
}

# Slide 1 title text (nested inside a group shape).
$oldTitle = "Prediction of Collision Cross-Section Values by Multimodal Graph Attention Network for Accurate Identification of Small Molecules"
$newTitle = "Accurate Prediction of Small Molecule Collision Cross-Section Values Through Chemical Class-Based Multimodal Graph Attention Network "

$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $top = $s1.Shapes.Item($i)
    if ($top.Type -eq 6) {
        # msoGroup - walk its items looking for the title text box.
        for ($j = 1; $j -le $top.GroupItems.Count; $j++) {
            $item = $top.GroupItems.Item($j)
            if ($item.HasTextFrame) {
                if ($item.TextFrame.TextRange.Text -eq $oldTitle) {
                    $item.TextFrame.TextRange.Text = $newTitle
                }
            }
        }
    } elseif ($top.HasTextFrame) {
        if ($top.TextFrame.TextRange.Text -eq $oldTitle) {
            $top.TextFrame.TextRange.Text = $newTitle
        }
    }
}
